$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 99.40000000000001
$ws.Cells.Item(5, 9).Value = 74.25
$ws.Cells.Item(5, 10).Value = 200
$ws.Cells.Item(5, 11).Value = 74.25
$ws.Cells.Item(5, 12).Value = 200
$ws.Cells.Item(5, 13).Value = 40.75
$ws.Cells.Item(5, 14).Value = -430
$ws.Cells.Item(40, 8).Value = 1500
$ws.Cells.Item(40, 10).Value = 1500
$ws.Cells.Item(40, 12).Value = 1500
$ws.Cells.Item(40, 14).Value = -1850
$ws.Cells.Item(94, 8).Value = 1139.8
$ws.Cells.Item(94, 9).Value = 924.75
$ws.Cells.Item(94, 11).Value = 924.75
$ws.Cells.Item(94, 13).Value = -473.75
$ws.Cells.Item(98, 8).Value = 971
$ws.Cells.Item(98, 9).Value = 874.381
$ws.Cells.Item(98, 10).Value = 3000
$ws.Cells.Item(98, 11).Value = 874.381
$ws.Cells.Item(98, 12).Value = 3000
$ws.Cells.Item(98, 13).Value = 623.619
$ws.Cells.Item(98, 14).Value = -5996
$ws.Cells.Item(101, 8).Value = 1269.5834
$ws.Cells.Item(101, 9).Value = 927.8570999999999
$ws.Cells.Item(101, 11).Value = 2783.5713
$ws.Cells.Item(101, 13).Value = -1161.5713
$ws.Cells.Item(113, 8).Value = 2550.3572
$ws.Cells.Item(113, 9).Value = 2543.5715
$ws.Cells.Item(113, 10).Value = 2557.1428
$ws.Cells.Item(113, 11).Value = 2543.5715
$ws.Cells.Item(113, 12).Value = 2557.1428
$ws.Cells.Item(113, 13).Value = 710.4285
$ws.Cells.Item(113, 14).Value = -9065.1428
$ws.Cells.Item(122, 8).Value = 971
$ws.Cells.Item(122, 9).Value = 874.381
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 2623.143
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -173.143
$ws.Cells.Item(122, 14).Value = -13900
$ws.Cells.Item(138, 8).Value = 4040.4578
$ws.Cells.Item(138, 9).Value = 4490.364
$ws.Cells.Item(138, 10).Value = 3971.7222
$ws.Cells.Item(138, 11).Value = 13471.092
$ws.Cells.Item(138, 12).Value = 11915.1666
$ws.Cells.Item(138, 13).Value = -8331.091999999999
$ws.Cells.Item(138, 14).Value = -22195.1666
$ws.Cells.Item(140, 8).Value = 71417.21000000001
$ws.Cells.Item(140, 10).Value = 71417.21000000001
$ws.Cells.Item(140, 12).Value = 71417.21000000001
$ws.Cells.Item(140, 14).Value = -81777.21000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 300
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).Value = $null
$ws.Cells.Item(32, 8).Value = 5756590
$ws.Cells.Item(32, 9).Value = 6418196.5
$ws.Cells.Item(32, 11).Value = 6418196.5
$ws.Cells.Item(32, 13).Value = -6417909.5
$ws.Cells.Item(63, 8).Value = 4159.7393
$ws.Cells.Item(63, 10).Value = 4791.154
$ws.Cells.Item(63, 12).Value = 4791.154
$ws.Cells.Item(63, 14).Value = -6163.154
$ws.Cells.Item(66, 8).Value = 4159.7393
$ws.Cells.Item(66, 10).Value = 4791.154
$ws.Cells.Item(66, 12).Value = 23955.77
$ws.Cells.Item(66, 14).Value = -30819.77

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 300
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).Value = $null
$ws.Cells.Item(22, 8).Value = 840
$ws.Cells.Item(22, 9).Value = 840
$ws.Cells.Item(22, 11).Value = 840
$ws.Cells.Item(22, 13).Value = -667
$ws.Cells.Item(68, 8).Value = 90000
$ws.Cells.Item(68, 10).Value = 90000
$ws.Cells.Item(68, 12).Value = 90000
$ws.Cells.Item(68, 14).Value = -91622
$ws.Cells.Item(71, 8).Value = 90000
$ws.Cells.Item(71, 10).Value = 90000
$ws.Cells.Item(71, 12).Value = 270000
$ws.Cells.Item(71, 14).Value = -278112
$ws.Cells.Item(100, 8).Value = 85000
$ws.Cells.Item(100, 10).Value = 85000
$ws.Cells.Item(100, 12).Value = 85000
$ws.Cells.Item(100, 14).Value = -87164
$ws.Cells.Item(105, 8).Value = 41668004
$ws.Cells.Item(105, 9).Value = 41668004
$ws.Cells.Item(105, 11).Value = 41668004
$ws.Cells.Item(105, 13).Value = -41666257

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 72.5
$ws.Cells.Item(7, 9).Value = 57.5
$ws.Cells.Item(7, 10).Value = 80
$ws.Cells.Item(7, 11).Value = 57.5
$ws.Cells.Item(7, 12).Value = 80
$ws.Cells.Item(7, 13).Value = 55.5
$ws.Cells.Item(7, 14).Value = -306
$ws.Cells.Item(31, 8).Value = 6186.239
$ws.Cells.Item(31, 9).Value = 1894.9231
$ws.Cells.Item(31, 11).Value = 1894.9231
$ws.Cells.Item(31, 13).Value = -1599.9231
$ws.Cells.Item(34, 8).Value = 6186.239
$ws.Cells.Item(34, 9).Value = 1894.9231
$ws.Cells.Item(34, 11).Value = 1894.9231
$ws.Cells.Item(34, 13).Value = -1692.9231
$ws.Cells.Item(58, 8).Value = 1657
$ws.Cells.Item(58, 9).Value = 1633.3334
$ws.Cells.Item(58, 10).Value = 1680.6666
$ws.Cells.Item(58, 11).Value = 1633.3334
$ws.Cells.Item(58, 12).Value = 1680.6666
$ws.Cells.Item(58, 13).Value = -1430.3334
$ws.Cells.Item(58, 14).Value = -2086.6666
$ws.Cells.Item(136, 8).Value = 1657
$ws.Cells.Item(136, 9).Value = 1633.3334
$ws.Cells.Item(136, 10).Value = 1680.6666
$ws.Cells.Item(136, 11).Value = 4900.0002
$ws.Cells.Item(136, 12).Value = 5041.9998
$ws.Cells.Item(136, 13).Value = -2350.0002
$ws.Cells.Item(136, 14).Value = -10141.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 6904.933
$ws.Cells.Item(68, 9).Value = 10661.333
$ws.Cells.Item(68, 10).Value = 1270.3334
$ws.Cells.Item(68, 11).Value = 31983.999
$ws.Cells.Item(68, 12).Value = 3811.0002
$ws.Cells.Item(68, 13).Value = -31172.999
$ws.Cells.Item(68, 14).Value = -5433.0002
$ws.Cells.Item(71, 8).Value = 6904.933
$ws.Cells.Item(71, 9).Value = 10661.333
$ws.Cells.Item(71, 10).Value = 1270.3334
$ws.Cells.Item(71, 11).Value = 95951.997
$ws.Cells.Item(71, 12).Value = 11433.0006
$ws.Cells.Item(71, 13).Value = -91895.997
$ws.Cells.Item(71, 14).Value = -19545.0006
$ws.Cells.Item(97, 8).Value = 13889901
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 13889901
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 41669703
$ws.Cells.Item(97, 13).Value = $null
$ws.Cells.Item(97, 14).Value = -41670695
$ws.Cells.Item(130, 8).Value = 2317.0715
$ws.Cells.Item(130, 9).Value = 1500
$ws.Cells.Item(130, 10).Value = 2379.923
$ws.Cells.Item(130, 11).Value = 4500
$ws.Cells.Item(130, 12).Value = 7139.768999999999
$ws.Cells.Item(130, 13).Value = 520
$ws.Cells.Item(130, 14).Value = -17179.769
$ws.Cells.Item(131, 8).Value = 4185.946
$ws.Cells.Item(131, 10).Value = 5163.1035
$ws.Cells.Item(131, 12).Value = 15489.3105
$ws.Cells.Item(131, 14).Value = -25569.3105
$ws.Cells.Item(132, 8).Value = 2835.25
$ws.Cells.Item(132, 9).Value = 2536.087
$ws.Cells.Item(132, 11).Value = 22824.783
$ws.Cells.Item(132, 13).Value = -20294.783
$ws.Cells.Item(137, 8).Value = 40753.035
$ws.Cells.Item(137, 10).Value = 86808.336
$ws.Cells.Item(137, 12).Value = 260425.008
$ws.Cells.Item(137, 14).Value = -270625.008
$ws.Cells.Item(138, 8).Value = 4445.4546
$ws.Cells.Item(138, 9).Value = 2237.5
$ws.Cells.Item(138, 10).Value = 10333.333
$ws.Cells.Item(138, 11).Value = 6712.5
$ws.Cells.Item(138, 12).Value = 30999.999
$ws.Cells.Item(138, 13).Value = -1572.5
$ws.Cells.Item(138, 14).Value = -41279.999
$ws.Cells.Item(140, 8).Value = 1840.3226
$ws.Cells.Item(140, 9).Value = 1320.8096
$ws.Cells.Item(140, 10).Value = 2931.3
$ws.Cells.Item(140, 11).Value = 3962.4288
$ws.Cells.Item(140, 12).Value = 8793.900000000001
$ws.Cells.Item(140, 13).Value = 1217.5712
$ws.Cells.Item(140, 14).Value = -19153.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5264.276
$ws.Cells.Item(70, 9).Value = 5098.2856
$ws.Cells.Item(70, 11).Value = 5098.2856
$ws.Cells.Item(70, 13).Value = -4828.2856
$ws.Cells.Item(73, 8).Value = 5264.276
$ws.Cells.Item(73, 9).Value = 5098.2856
$ws.Cells.Item(73, 11).Value = 5098.2856
$ws.Cells.Item(73, 13).Value = -4162.2856
$ws.Cells.Item(113, 8).Value = 86688.16
$ws.Cells.Item(113, 9).Value = 111923.3
$ws.Cells.Item(113, 11).Value = 111923.3
$ws.Cells.Item(113, 13).Value = -109753.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9332.111000000001
$ws.Cells.Item(7, 9).Value = 8797.333000000001
$ws.Cells.Item(7, 11).Value = 8797.333000000001
$ws.Cells.Item(7, 13).Value = -8685.333000000001
$ws.Cells.Item(40, 8).Value = 3530.2144
$ws.Cells.Item(40, 9).Value = 3301.8
$ws.Cells.Item(40, 10).Value = 4101.25
$ws.Cells.Item(40, 11).Value = 3301.8
$ws.Cells.Item(40, 12).Value = 4101.25
$ws.Cells.Item(40, 13).Value = -3165.8
$ws.Cells.Item(40, 14).Value = -4373.25
$ws.Cells.Item(126, 8).Value = 9332.111000000001
$ws.Cells.Item(126, 9).Value = 8797.333000000001
$ws.Cells.Item(126, 11).Value = 26391.999
$ws.Cells.Item(126, 13).Value = -23921.999
$ws.Cells.Item(132, 8).Value = 4272
$ws.Cells.Item(132, 9).Value = 3199.6
$ws.Cells.Item(132, 10).Value = 5165.6665
$ws.Cells.Item(132, 11).Value = 9598.799999999999
$ws.Cells.Item(132, 12).Value = 15496.9995
$ws.Cells.Item(132, 13).Value = -7068.799999999999
$ws.Cells.Item(132, 14).Value = -20556.9995
$ws.Cells.Item(136, 8).Value = 2105.5
$ws.Cells.Item(136, 9).Value = 2005.7894
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 6017.3682
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -3467.3682
$ws.Cells.Item(136, 14).Value = -17100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1130.2858
$ws.Cells.Item(113, 9).Value = 1356.7273
$ws.Cells.Item(113, 11).Value = 4070.1819
$ws.Cells.Item(113, 13).Value = -1900.1819
$ws.Cells.Item(116, 8).Value = 69170
$ws.Cells.Item(116, 10).Value = 69170
$ws.Cells.Item(116, 12).Value = 69170
$ws.Cells.Item(116, 14).Value = -78348
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 2518.0938
$ws.Cells.Item(122, 9).Value = 1733.95
$ws.Cells.Item(122, 10).Value = 3825
$ws.Cells.Item(122, 11).Value = 5201.85
$ws.Cells.Item(122, 12).Value = 11475
$ws.Cells.Item(122, 13).Value = -2751.85
$ws.Cells.Item(122, 14).Value = -16375
$ws.Cells.Item(126, 8).Value = 1568.9048
$ws.Cells.Item(126, 9).Value = 1652.9375
$ws.Cells.Item(126, 10).Value = 1300
$ws.Cells.Item(126, 11).Value = 4958.8125
$ws.Cells.Item(126, 12).Value = 3900
$ws.Cells.Item(126, 13).Value = -2488.8125
$ws.Cells.Item(126, 14).Value = -8840
$ws.Cells.Item(132, 8).Value = 9410948
$ws.Cells.Item(132, 9).Value = 2296.6155
$ws.Cells.Item(132, 11).Value = 6889.8465
$ws.Cells.Item(132, 13).Value = -4359.8465
